$d = $word.ActiveDocument

# Helper: replace the text of a ListBullet paragraph ("What we like" / "What we
# don't like" items) while preserving its exact original run layout (a leading
# empty <w:r/> run followed by the text run). Plain Find/Replace on these
# paragraphs collapses that leading empty run away because the matched run has
# no distinguishing run properties, so we rebuild the paragraph from a minimal
# OOXML fragment via Range.InsertXML instead.
function Set-BulletParagraphText($paragraphIndex, $newText) {
    $para = $d.Paragraphs($paragraphIndex)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body>' `
        + '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' `
        + '<w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p>' `
        + '</w:body></w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

# 1. Title/heading. This exact phrase occurs twice (the Heading1 at the top
#    of the document, and a bold recap run near the bottom) and both should
#    get the same new text, so a single "Replace All" covers both spots.
$d.Content.Find.Execute("Play Bonsai Spins Free Slot Machine Online", $true, $false, $false, $false, $false, $true, 1, $false, "Play Bonsai Spins for Free - Exciting Gameplay and Beautiful Design", 2)

# 2. "What we like" bullet list.
Set-BulletParagraphText 42 "Exciting gameplay mechanism with Wild and Scatter symbols"
Set-BulletParagraphText 43 "Outstanding graphics and well-executed design"
Set-BulletParagraphText 44 "Engaging and relaxing Oriental music soundtrack"
Set-BulletParagraphText 45 "Very accessible and user-friendly"

# 3. "What we don't like" bullet list.
Set-BulletParagraphText 47 "None of note"
Set-BulletParagraphText 48 "No progressive jackpot"

# 4. Italic meta description at the very end.
$d.Content.Find.Execute("Discover the exciting Bonsai Spins slot machine game with Wild and Scatter symbols. Play now for free and experience beautiful graphics and Oriental music.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Bonsai Spins and play for free. Exciting gameplay, stunning graphics, and relaxing music.", 2)
